# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) counts across the four sheets,
# and marks one "最低票价" (G column) entry as "不可售" (not sellable).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Changes = @{
    "F2"  = 754
    "F4"  = 2042
    "F5"  = 6114
    "F6"  = 3473
    "F9"  = 1448
    "F10" = 4821
    "F11" = 1116
    "F12" = 1797
    "F13" = 26
    "F16" = 221
    "F18" = 1066
    "F19" = 339
    "F21" = 39
    "F23" = 7
    "F24" = 224
    "F25" = 112
    "F26" = 31
    "F27" = 1158
    "F28" = 436
    "F30" = 243
    "F31" = 502
    "F32" = 1003
    "F33" = 32
    "F34" = 1862
    "F35" = 2318
    "F36" = 1092
    "F38" = 28
    "F40" = 103
    "F41" = 677
    "F42" = 505
    "F43" = 69
    "F44" = 697
    "F45" = 61
    "F46" = 481
    "F47" = 506
    "F48" = 247
    "F49" = 164
}
foreach ($addr in $sheet1Changes.Keys) {
    $ws1.Range($addr).Value = $sheet1Changes[$addr]
}

# ---- Sheet "演出" (Performance) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "不可售"
$ws2.Range("F17").Value = 95
$ws2.Range("F26").Value = 32

# ---- Sheet "本地生活" (Local Life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 832

# ---- Sheet "全部类型" (All Types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Changes = @{
    "F2"  = 832
    "F3"  = 754
    "F5"  = 2042
    "F6"  = 3473
    "F8"  = 1448
    "F9"  = 4821
    "F10" = 1797
    "F11" = 26
    "F17" = 221
    "F20" = 1066
    "F21" = 339
    "F23" = 224
    "F25" = 31
    "F26" = 1158
    "F27" = 436
    "F29" = 243
    "F31" = 1003
    "F32" = 1862
    "F33" = 2318
    "F34" = 1092
    "F36" = 28
    "F38" = 103
    "F41" = 677
    "F42" = 505
    "F43" = 697
    "F44" = 481
    "F45" = 506
    "F46" = 247
    "F48" = 164
    "F49" = 32
}
foreach ($addr in $sheet4Changes.Keys) {
    $ws4.Range($addr).Value = $sheet4Changes[$addr]
}

$wb.Save()
